$d = $word.ActiveDocument

$replacements = @(
    @{Old="334÷8="; New="383÷3="},
    @{Old="809÷2="; New="658÷2="},
    @{Old="663÷2="; New="724÷8="},
    @{Old="595÷5="; New="399÷6="},
    @{Old="368÷2="; New="611÷9="},
    @{Old="823÷3="; New="465÷8="},
    @{Old="295÷8="; New="591÷3="},
    @{Old="152÷9="; New="386÷2="},
    @{Old="108÷7="; New="228÷4="},
    @{Old="434÷8="; New="609÷6="},
    @{Old="973÷9="; New="543÷8="},
    @{Old="702÷8="; New="582÷2="},
    @{Old="912÷6="; New="419÷7="},
    @{Old="597÷5="; New="985÷2="},
    @{Old="646÷6="; New="242÷2="},
    @{Old="297÷9="; New="940÷8="},
    @{Old="999÷7="; New="514÷7="},
    @{Old="847÷3="; New="360÷6="},
    @{Old="267÷3="; New="474÷6="},
    @{Old="201÷5="; New="284÷6="},
    @{Old="254÷5="; New="855÷6="},
    @{Old="264÷3="; New="397÷7="},
    @{Old="632÷2="; New="631÷4="},
    @{Old="763÷8="; New="110÷8="},
    @{Old="969÷9="; New="498÷6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
